# Create biodiversity indicators for age and species
# Update the emission results with recalculated values for solid_wood,
# sum_product, ecosystem and system columns (B, D, E, F) across rows 2-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  B = 0;                  E = -2515007.418206576 },
    @{ Row = 3;  B = 876663.5182445209;  E = 1306319.689602927 },
    @{ Row = 4;  B = 1165700.198931394;  E = -401483.2380664512 },
    @{ Row = 5;  B = 1390816.352275759;  E = -1502979.778848212 },
    @{ Row = 6;  B = 1453091.396268783;  E = -1473402.664590993 },
    @{ Row = 7;  B = 1485408.60282823;   E = -759328.4909590487 },
    @{ Row = 8;  B = 1417358.485838958;  E = -2833412.93896434 },
    @{ Row = 9;  B = 1486404.488605577;  E = -1444270.181108203 },
    @{ Row = 10; B = 1237228.929894109;  E = -3076644.03439256 },
    @{ Row = 11; B = 979493.4069906169;  E = -7738433.732811249 }
)

foreach ($item in $data) {
    $r = [int]$item.Row
    $bVal = [double]$item.B
    $eVal = [double]$item.E
    $cVal = [double]($ws.Cells.Item($r, 3).Value2)
    $dVal = $bVal + $cVal
    $fVal = $dVal + $eVal

    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 4).Value = $dVal
    $ws.Cells.Item($r, 5).Value = $eVal
    $ws.Cells.Item($r, 6).Value = $fVal
}
